# This edit swaps the data (Code, Item, Cost, Price, Qty, Total columns B:G)
# between specific pairs of adjacent rows in the stock report. The row's
# serial number in column A is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("B$r1" + ":G$r1")
    $range2 = $ws.Range("B$r2" + ":G$r2")
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2
    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-Rows 151 152
Swap-Rows 198 199
Swap-Rows 228 229
Swap-Rows 237 238
Swap-Rows 326 327
Swap-Rows 371 372
Swap-Rows 373 374
Swap-Rows 401 402
Swap-Rows 554 555
Swap-Rows 563 564
Swap-Rows 568 569
Swap-Rows 573 574
Swap-Rows 644 645
Swap-Rows 663 664
Swap-Rows 673 674
Swap-Rows 834 835
